$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Remove the INVOICE column - COST shifts left to become column B
$ws.Columns(2).Delete()

# Update the order row with the new import data from cbip
$ws.Range("A2").Value = 394732
$ws.Range("B2").Value = 20

# Restyle the ID cell with the new (Inter) font used for imported rows
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Font.Name = "Inter"
$ws.Range("A2").Font.Color = 5526612

$ws.Range("C11").Select()
